# Actualizacion de los test
# Reproduces the data "rotation" in Hoja1: a new method column (Msucesiva)
# moves from the last position (H / Q) to the first position (B / K) in the
# two result tables (rows 20-23), row 15 gets the missing "n"/"n2" labels
# filled in for the right-hand table (K15:Q15), and the view is scrolled
# back to the top with a changed zoom/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 15: fill in the complexity labels for the right-hand table
# (n = "n", n2 = "n2") that were left blank before.
# ---------------------------------------------------------------------
$ws.Range("K15").Value = "n"
$ws.Range("L15").Value = "n2"
$ws.Range("M15").Value = "n2"
$ws.Range("N15").Value = "n2"
$ws.Range("O15").Value = "n"
$ws.Range("Q15").Value = "n"

# ---------------------------------------------------------------------
# Helper: rotate the 7 values held in columns B..H (or K..Q) of a given
# row one position to the right, wrapping the last column's value back
# to the first column ($null means "leave the cell empty").
# ---------------------------------------------------------------------
function Rotate-Row($ws, $row, $cols) {
    $n = $cols.Length
    $vals = @()
    foreach ($c in $cols) {
        $cell = $ws.Range($c + $row)
        if ($cell.Value2 -eq $null) {
            $vals += ,$null
        } else {
            $vals += ,$cell.Value2
        }
    }
    for ($i = 0; $i -lt $n; $i++) {
        $srcIdx = ($i - 1 + $n) % $n
        $destCol = $cols[$i]
        $val = $vals[$srcIdx]
        $cell = $ws.Range($destCol + $row)
        if ($val -eq $null) {
            $cell.ClearContents()
        } else {
            $cell.Value = $val
        }
    }
}

$leftCols = @("B", "C", "D", "E", "F", "G", "H")
$rightCols = @("K", "L", "M", "N", "O", "P", "Q")

# Row 20 holds the header labels (method names) for both tables.
Rotate-Row $ws 20 $leftCols
Rotate-Row $ws 20 $rightCols

# Rows 21-23 hold the measured values for both tables.
Rotate-Row $ws 21 $leftCols
Rotate-Row $ws 21 $rightCols
Rotate-Row $ws 22 $leftCols
Rotate-Row $ws 22 $rightCols
Rotate-Row $ws 23 $leftCols
Rotate-Row $ws 23 $rightCols

# ---------------------------------------------------------------------
# Sheet view: scroll back to the top-left corner, zoom to 115%, and move
# the active selection to F18.
# ---------------------------------------------------------------------
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$win.Zoom = 115
$ws.Range("F18").Select()
